$wb = $excel.ActiveWorkbook

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 3792.6206
$ws.Cells.Item(64, 9).Value = 3647.0588
$ws.Cells.Item(64, 10).Value = 3998.8333
$ws.Cells.Item(64, 11).Value = 3647.0588
$ws.Cells.Item(64, 12).Value = 3998.8333
$ws.Cells.Item(64, 13).Value = -3399.0588
$ws.Cells.Item(64, 14).Value = -4494.8333

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(67, 8).Value = 3792.6206
$ws.Cells.Item(67, 9).Value = 3647.0588
$ws.Cells.Item(67, 10).Value = 3998.8333
$ws.Cells.Item(67, 11).Value = 3647.0588
$ws.Cells.Item(67, 12).Value = 3998.8333
$ws.Cells.Item(67, 13).Value = -2789.0588
$ws.Cells.Item(67, 14).Value = -5714.8333

# ALC row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 10982012
$ws.Cells.Item(70, 9).Value = 13828711
$ws.Cells.Item(70, 10).Value = 1885.5714
$ws.Cells.Item(70, 11).Value = 41486133
$ws.Cells.Item(70, 12).Value = 5656.7142
$ws.Cells.Item(70, 13).Value = -41485863
$ws.Cells.Item(70, 14).Value = -6196.7142

# ALC row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value = 10982012
$ws.Cells.Item(73, 9).Value = 13828711
$ws.Cells.Item(73, 10).Value = 1885.5714
$ws.Cells.Item(73, 11).Value = 41486133
$ws.Cells.Item(73, 12).Value = 5656.7142
$ws.Cells.Item(73, 13).Value = -41485197
$ws.Cells.Item(73, 14).Value = -7528.7142

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 4526.5454
$ws.Cells.Item(76, 9).Value = 3370.2856
$ws.Cells.Item(76, 11).Value = 3370.2856
$ws.Cells.Item(76, 13).Value = -3055.2856

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(79, 8).Value = 4526.5454
$ws.Cells.Item(79, 9).Value = 3370.2856
$ws.Cells.Item(79, 11).Value = 3370.2856
$ws.Cells.Item(79, 13).Value = -2278.2856

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 38201.63
$ws.Cells.Item(100, 9).Value = 63558.75
$ws.Cells.Item(100, 10).Value = 1318.5454
$ws.Cells.Item(100, 11).Value = 63558.75
$ws.Cells.Item(100, 12).Value = 1318.5454
$ws.Cells.Item(100, 13).Value = -63017.75
$ws.Cells.Item(100, 14).Value = -2400.5454

# ALC row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(103, 8).Value = 1234.5454
$ws.Cells.Item(103, 9).Value = 813.3333
$ws.Cells.Item(103, 10).Value = 1740
$ws.Cells.Item(103, 11).Value = 2439.9999
$ws.Cells.Item(103, 12).Value = 5220
$ws.Cells.Item(103, 13).Value = -1853.9999
$ws.Cells.Item(103, 14).Value = -6392

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 485967.75
$ws.Cells.Item(32, 9).Value = 2027.3658
$ws.Cells.Item(32, 10).Value = 1430803.8
$ws.Cells.Item(32, 11).Value = 2027.3658
$ws.Cells.Item(32, 12).Value = 1430803.8
$ws.Cells.Item(32, 13).Value = -1740.3658
$ws.Cells.Item(32, 14).Value = -1431377.8

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 5566.3022
$ws.Cells.Item(61, 9).Value = 6082.2905
$ws.Cells.Item(61, 10).Value = 4233.3335
$ws.Cells.Item(61, 11).Value = 6082.2905
$ws.Cells.Item(61, 12).Value = 4233.3335
$ws.Cells.Item(61, 13).Value = -5870.2905
$ws.Cells.Item(61, 14).Value = -4657.3335

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 3708192.5
$ws.Cells.Item(63, 9).Value = 16667916
$ws.Cells.Item(63, 10).Value = 5414.2856
$ws.Cells.Item(63, 11).Value = 16667916
$ws.Cells.Item(63, 12).Value = 5414.2856
$ws.Cells.Item(63, 13).Value = -16667230
$ws.Cells.Item(63, 14).Value = -6786.2856

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value = 3708192.5
$ws.Cells.Item(66, 9).Value = 16667916
$ws.Cells.Item(66, 10).Value = 5414.2856
$ws.Cells.Item(66, 11).Value = 83339580
$ws.Cells.Item(66, 12).Value = 27071.428
$ws.Cells.Item(66, 13).Value = -83336148
$ws.Cells.Item(66, 14).Value = -33935.428

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 876.8
$ws.Cells.Item(97, 9).Value = 301.92307
$ws.Cells.Item(97, 10).Value = 1944.4286
$ws.Cells.Item(97, 11).Value = 301.92307
$ws.Cells.Item(97, 12).Value = 1944.4286
$ws.Cells.Item(97, 13).Value = 194.07693
$ws.Cells.Item(97, 14).Value = -2936.4286

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 2529.75
$ws.Cells.Item(102, 9).Value = 2848
$ws.Cells.Item(102, 10).Value = 1999.3334
$ws.Cells.Item(102, 11).Value = 2848
$ws.Cells.Item(102, 12).Value = 1999.3334
$ws.Cells.Item(102, 13).Value = -1226
$ws.Cells.Item(102, 14).Value = -5243.3334

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 5566.3022
$ws.Cells.Item(136, 9).Value = 6082.2905
$ws.Cells.Item(136, 10).Value = 4233.3335
$ws.Cells.Item(136, 11).Value = 18246.8715
$ws.Cells.Item(136, 12).Value = 12700.0005
$ws.Cells.Item(136, 13).Value = -15696.8715
$ws.Cells.Item(136, 14).Value = -17800.0005

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 667.1111
$ws.Cells.Item(94, 9).Value = 600.0952
$ws.Cells.Item(94, 10).Value = 901.6667
$ws.Cells.Item(94, 11).Value = 600.0952
$ws.Cells.Item(94, 12).Value = 901.6667
$ws.Cells.Item(94, 13).Value = -149.0952
$ws.Cells.Item(94, 14).Value = -1803.6667

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1089.75
$ws.Cells.Item(99, 9).Value = 708.41174
$ws.Cells.Item(99, 10).Value = 1679.091
$ws.Cells.Item(99, 11).Value = 708.41174
$ws.Cells.Item(99, 12).Value = 1679.091
$ws.Cells.Item(99, 13).Value = 789.58826
$ws.Cells.Item(99, 14).Value = -4675.091

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2283.889
$ws.Cells.Item(105, 9).Value = 2380
$ws.Cells.Item(105, 10).Value = 2163.75
$ws.Cells.Item(105, 11).Value = 2380
$ws.Cells.Item(105, 12).Value = 2163.75
$ws.Cells.Item(105, 13).Value = -633
$ws.Cells.Item(105, 14).Value = -5657.75

# CRP row 18
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(18, 8).Value = 23000
$ws.Cells.Item(18, 10).Value = 23000
$ws.Cells.Item(18, 12).Value = 23000
$ws.Cells.Item(18, 14).Value = -23460

# CUL row 80
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 3272.4827
$ws.Cells.Item(80, 9).Value = 1267.3334
$ws.Cells.Item(80, 10).Value = 3503.8462
$ws.Cells.Item(80, 11).Value = 3802.0002
$ws.Cells.Item(80, 12).Value = 10511.5386
$ws.Cells.Item(80, 13).Value = -2866.0002
$ws.Cells.Item(80, 14).Value = -12383.5386

# CUL row 83
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(83, 8).Value = 3272.4827
$ws.Cells.Item(83, 9).Value = 1267.3334
$ws.Cells.Item(83, 10).Value = 3503.8462
$ws.Cells.Item(83, 11).Value = 11406.0006
$ws.Cells.Item(83, 12).Value = 31534.6158
$ws.Cells.Item(83, 13).Value = -6726.000599999999
$ws.Cells.Item(83, 14).Value = -40894.6158

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 40082.605
$ws.Cells.Item(70, 9).Value = 70807.2
$ws.Cells.Item(70, 10).Value = 4631.154
$ws.Cells.Item(70, 11).Value = 70807.2
$ws.Cells.Item(70, 12).Value = 4631.154
$ws.Cells.Item(70, 13).Value = -70537.2
$ws.Cells.Item(70, 14).Value = -5171.154

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 40082.605
$ws.Cells.Item(73, 9).Value = 70807.2
$ws.Cells.Item(73, 10).Value = 4631.154
$ws.Cells.Item(73, 11).Value = 70807.2
$ws.Cells.Item(73, 12).Value = 4631.154
$ws.Cells.Item(73, 13).Value = -69871.2
$ws.Cells.Item(73, 14).Value = -6503.154

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3177
$ws.Cells.Item(80, 9).Value = 2752.5
$ws.Cells.Item(80, 10).Value = 3460
$ws.Cells.Item(80, 11).Value = 2752.5
$ws.Cells.Item(80, 12).Value = 3460
$ws.Cells.Item(80, 13).Value = -1754.5
$ws.Cells.Item(80, 14).Value = -5456

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 3177
$ws.Cells.Item(83, 9).Value = 2752.5
$ws.Cells.Item(83, 10).Value = 3460
$ws.Cells.Item(83, 11).Value = 13762.5
$ws.Cells.Item(83, 12).Value = 17300
$ws.Cells.Item(83, 13).Value = -8770.5
$ws.Cells.Item(83, 14).Value = -27284

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1406.55
$ws.Cells.Item(97, 9).Value = 1390
$ws.Cells.Item(97, 10).Value = 1437.2858
$ws.Cells.Item(97, 11).Value = 1390
$ws.Cells.Item(97, 12).Value = 1437.2858
$ws.Cells.Item(97, 13).Value = -894
$ws.Cells.Item(97, 14).Value = -2429.2858

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 71832.92999999999
$ws.Cells.Item(107, 9).Value = 111486.78
$ws.Cells.Item(107, 10).Value = 456
$ws.Cells.Item(107, 11).Value = 111486.78
$ws.Cells.Item(107, 12).Value = 456
$ws.Cells.Item(107, 13).Value = -109566.78
$ws.Cells.Item(107, 14).Value = -4296

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 1647.303
$ws.Cells.Item(68, 9).Value = 1465.2593
$ws.Cells.Item(68, 10).Value = 2466.5
$ws.Cells.Item(68, 11).Value = 1465.2593
$ws.Cells.Item(68, 12).Value = 2466.5
$ws.Cells.Item(68, 13).Value = -716.2592999999999
$ws.Cells.Item(68, 14).Value = -3964.5

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value = 1647.303
$ws.Cells.Item(71, 9).Value = 1465.2593
$ws.Cells.Item(71, 10).Value = 2466.5
$ws.Cells.Item(71, 11).Value = 7326.296499999999
$ws.Cells.Item(71, 12).Value = 12332.5
$ws.Cells.Item(71, 13).Value = -3582.296499999999
$ws.Cells.Item(71, 14).Value = -19820.5

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 2346.3157
$ws.Cells.Item(82, 9).Value = 1385.6364
$ws.Cells.Item(82, 10).Value = 3667.25
$ws.Cells.Item(82, 11).Value = 1385.6364
$ws.Cells.Item(82, 12).Value = 3667.25
$ws.Cells.Item(82, 13).Value = -1024.6364
$ws.Cells.Item(82, 14).Value = -4389.25

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(85, 8).Value = 2346.3157
$ws.Cells.Item(85, 9).Value = 1385.6364
$ws.Cells.Item(85, 10).Value = 3667.25
$ws.Cells.Item(85, 11).Value = 1385.6364
$ws.Cells.Item(85, 12).Value = 3667.25
$ws.Cells.Item(85, 13).Value = -137.6364000000001
$ws.Cells.Item(85, 14).Value = -6163.25

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 62031.53
$ws.Cells.Item(93, 9).Value = 1707.8182
$ws.Cells.Item(93, 10).Value = 172625
$ws.Cells.Item(93, 11).Value = 1707.8182
$ws.Cells.Item(93, 12).Value = 172625
$ws.Cells.Item(93, 13).Value = -459.8181999999999
$ws.Cells.Item(93, 14).Value = -175121

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 1937.7
$ws.Cells.Item(100, 9).Value = 1414.2858
$ws.Cells.Item(100, 10).Value = 2219.5386
$ws.Cells.Item(100, 11).Value = 1414.2858
$ws.Cells.Item(100, 12).Value = 2219.5386
$ws.Cells.Item(100, 13).Value = -873.2858000000001
$ws.Cells.Item(100, 14).Value = -3301.5386

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 6820824
$ws.Cells.Item(132, 9).Value = 8335257
$ws.Cells.Item(132, 10).Value = 3791957.8
$ws.Cells.Item(132, 11).Value = 25005771
$ws.Cells.Item(132, 12).Value = 11375873.4
$ws.Cells.Item(132, 13).Value = -25003241
$ws.Cells.Item(132, 14).Value = -11380933.4

# WVR row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 5325.8335
$ws.Cells.Item(96, 9).Value = 9418.333000000001
$ws.Cells.Item(96, 10).Value = 1233.3334
$ws.Cells.Item(96, 11).Value = 9418.333000000001
$ws.Cells.Item(96, 12).Value = 1233.3334
$ws.Cells.Item(96, 13).Value = -8045.333000000001
$ws.Cells.Item(96, 14).Value = -3979.3334

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 338.41666
$ws.Cells.Item(100, 9).Value = 350.8889
$ws.Cells.Item(100, 10).Value = 301
$ws.Cells.Item(100, 11).Value = 701.7778
$ws.Cells.Item(100, 12).Value = 602
$ws.Cells.Item(100, 13).Value = -160.7778
$ws.Cells.Item(100, 14).Value = -1684

Write-Host "Applied all Fenrir_Profits market price updates"
